$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Chris Paul"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "San Antonio Spurs"
$ws.Range("A3").Value = "Russell Westbrook"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Denver Nuggets"
$ws.Range("A4").Value = "Payton Pritchard"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Boston Celtics"
$ws.Range("A5").Value = "Nick Smith Jr."
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Charlotte Hornets"
$ws.Range("A6").Value = "Ayo Dosunmu"
$ws.Range("B6").Value = "PG,SG,SF"
$ws.Range("C6").Value = "Chicago Bulls"
$ws.Range("A7").Value = "Paolo Banchero"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Orlando Magic"
$ws.Range("A8").Value = "Jaylen Brown"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Boston Celtics"
$ws.Range("A9").Value = "Pascal Siakam"
$ws.Range("B9").Value = "SF,PF,C"
$ws.Range("C9").Value = "Indiana Pacers"
$ws.Range("A10").Value = "Nikola Jokic"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Denver Nuggets"
$ws.Range("A11").Value = "Rudy Gobert"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Minnesota Timberwolves"
$ws.Range("A12").Value = "Dejounte Murray"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "New Orleans Pelicans"
$ws.Range("A13").Value = "Deni Avdija"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Portland Trail Blazers"
$ws.Range("A14").Value = "Jalen Green"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Houston Rockets"
$ws.Range("A15").Value = "Jose Alvarado"
$ws.Range("B15").Value = "PG"
$ws.Range("C15").Value = "New Orleans Pelicans"
$ws.Range("A16").Value = "Jakob Poeltl"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Toronto Raptors"
$ws.Range("A17").Value = "Chet Holmgren"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Oklahoma City Thunder"
$ws.Range("A18").Value = "Jalen Suggs"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Orlando Magic"
$ws.Range("A19").Value = "Khris Middleton"
$ws.Range("B19").Value = "SF"
$ws.Range("C19").Value = "Milwaukee Bucks"
